$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "pyqelbdg@nldteery.com"
$ws.Range("B1").Value = "kJOixCjbu1Y"
$ws.Range("C1").Value = "ewduymhkxxz"

$ws.Range("A5").Select()
